{"js": "const pairs = [\n  [\"2023-12-26 Tuesday\", \"2023-12-27 Wednesday\"],\n  [\"33\u00f75=6, 3\", \"27\u00f78=3, 3\"],\n  [\"57\u00f73=19, 0\", \"93\u00f79=10, 3\"],\n  [\"66\u00f72=33, 0\", \"30\u00f73=10, 0\"],\n  [\"11\u00f75=2, 1\", \"42\u00f73=14, 0\"],\n  [\"23\u00f77=3, 2\", \"44\u00f74=11, 0\"],\n  [\"72\u00f72=36, 0\", \"72\u00f75=14, 2\"],\n  [\"18\u00f72=9, 0\", \"34\u00f78=4, 2\"],\n  [\"35\u00f74=8, 3\", \"96\u00f79=10, 6\"],\n  [\"31\u00f72=15, 1\", \"83\u00f73=27, 2\"],\n  [\"29\u00f76=4, 5\", \"90\u00f75=18, 0\"],\n  [\"25\u00f76=4, 1\", \"51\u00f76=8, 3\"],\n  [\"96\u00f73=32, 0\", \"92\u00f75=18, 2\"],\n  [\"97\u00f73=32, 1\", \"65\u00f74=16, 1\"],\n  [\"38\u00f72=19, 0\", \"19\u00f75=3, 4\"],\n  [\"86\u00f77=12, 2\", \"87\u00f76=14, 3\"],\n  [\"85\u00f74=21, 1\", \"30\u00f73=10, 0\"],\n  [\"27\u00f74=6, 3\", \"38\u00f78=4, 6\"],\n  [\"83\u00f78=10, 3\", \"38\u00f77=5, 3\"],\n  [\"74\u00f73=24, 2\", \"72\u00f79=8, 0\"],\n  [\"12\u00f79=1, 3\", \"52\u00f75=10, 2\"],\n  [\"24\u00f78=3, 0\", \"25\u00f78=3, 1\"],\n  [\"39\u00f76=6, 3\", \"73\u00f77=10, 3\"],\n  [\"59\u00f75=11, 4\", \"50\u00f78=6, 2\"],\n  [\"52\u00f77=7, 3\", \"37\u00f77=5, 2\"],\n  [\"65\u00f72=32, 1\", \"87\u00f74=21, 3\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@('2023-12-26 Tuesday', '2023-12-27 Wednesday')\n    ,@('33\u00f75=6, 3', '27\u00f78=3, 3')\n    ,@('57\u00f73=19, 0', '93\u00f79=10, 3')\n    ,@('66\u00f72=33, 0', '30\u00f73=10, 0')\n    ,@('11\u00f75=2, 1', '42\u00f73=14, 0')\n    ,@('23\u00f77=3, 2', '44\u00f74=11, 0')\n    ,@('72\u00f72=36, 0', '72\u00f75=14, 2')\n    ,@('18\u00f72=9, 0', '34\u00f78=4, 2')\n    ,@('35\u00f74=8, 3', '96\u00f79=10, 6')\n    ,@('31\u00f72=15, 1', '83\u00f73=27, 2')\n    ,@('29\u00f76=4, 5', '90\u00f75=18, 0')\n    ,@('25\u00f76=4, 1', '51\u00f76=8, 3')\n    ,@('96\u00f73=32, 0', '92\u00f75=18, 2')\n    ,@('97\u00f73=32, 1', '65\u00f74=16, 1')\n    ,@('38\u00f72=19, 0', '19\u00f75=3, 4')\n    ,@('86\u00f77=12, 2', '87\u00f76=14, 3')\n    ,@('85\u00f74=21, 1', '30\u00f73=10, 0')\n    ,@('27\u00f74=6, 3', '38\u00f78=4, 6')\n    ,@('83\u00f78=10, 3', '38\u00f77=5, 3')\n    ,@('74\u00f73=24, 2', '72\u00f79=8, 0')\n    ,@('12\u00f79=1, 3', '52\u00f75=10, 2')\n    ,@('24\u00f78=3, 0', '25\u00f78=3, 1')\n    ,@('39\u00f76=6, 3', '73\u00f77=10, 3')\n    ,@('59\u00f75=11, 4', '50\u00f78=6, 2')\n    ,@('52\u00f77=7, 3', '37\u00f77=5, 2')\n    ,@('65\u00f72=32, 1', '87\u00f74=21, 3')\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll) | Out-Null\n}"}
